$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.673.83"
$ws.Range("E2").Value = "  -3.08%  "
$ws.Range("D3").Value = "1.851.72"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.61"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4661"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3910"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.41"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07910"
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9841"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.37"
$ws.Range("E12").Value = "  -5.54%  "
$ws.Range("D13").Value = "1.817.03"
$ws.Range("E13").Value = "  -6.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.850"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.999"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06843"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.68"
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001008"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.12"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -1.00%  "
$ws.Range("D22").Value = "28.691.53"
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.390"
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  -5.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.137"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "2.083.72"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.36"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.48"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.101"
$ws.Range("E29").Value = "  -5.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.024"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.60"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9777"
$ws.Range("E32").Value = "  -3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09437"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.371"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.352"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06174"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02201"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.164"
$ws.Range("E39").Value = "  -1.69%  "
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.604"
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.21"
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.370"
$ws.Range("E44").Value = "  -3.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.250"
$ws.Range("E45").Value = "  -2.78%  "
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.77"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07163"
$ws.Range("E48").Value = "  -4.52%  "
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.23"
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.33"
$ws.Range("E51").Value = "  +3.25%  "
